# Natmi LR-pair output "Inhbb-Acvr2a" was regenerated following Dr Hou's
# advice: the cell-cluster set now includes "FAPs" as a sending cluster too
# (in addition to the existing "ECs" and "sCs"), which triples the number
# of Sending-cluster x Target-cluster combinations (3 x 4 = 12 rows instead
# of 2 x 4 = 8) and refreshes every computed statistic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: SendingCluster, Ligand, Receptor, TargetCluster, then the 16
# numeric columns E..T in order.
$rows = @(
    @("ECs","Inhbb","Acvr2a","ECs", 3,1,2.525153666666667,7.575461000000001,0.7378191762484796,0.7378191762484798,3,1,17.13024733333333,51.390742,0.2959211466465044,0.2959211466465043,43.25650686467355,389.308561782062,0.2183362966532294,0.2183362966532294),
    @("ECs","Inhbb","Acvr2a","FAPs",3,1,2.525153666666667,7.575461000000001,0.7378191762484796,0.7378191762484798,3,1,27.61090666666666,82.83271999999999,0.4769721651858779,0.4769721651858778,69.72178220932445,627.49603988392,0.3519192100108982,0.3519192100108982),
    @("ECs","Inhbb","Acvr2a","M2",  3,1,2.525153666666667,7.575461000000001,0.7378191762484796,0.7378191762484798,3,1,4.423514,13.270542,0.07641520344774541,0.0764152034477454,11.17005259665133,100.530473369862,0.0563806024606755,0.05638060246067551),
    @("ECs","Inhbb","Acvr2a","sCs", 3,1,2.525153666666667,7.575461000000001,0.7378191762484796,0.7378191762484798,3,1,8.723210333333332,26.169631,0.1506914847198724,0.1506914847198724,22.02744655832122,198.247019024891,0.1111830671236766,0.1111830671236766),

    @("FAPs","Inhbb","Acvr2a","ECs", 2,0.6666666666666666,0.3739756666666667,1.121927,0.1092711394000877,0.1092711394000878,3,1,17.13024733333333,51.390742,0.2959211466465044,0.2959211466465043,6.406295666648223,57.656660999834,0.03233564086664398,0.03233564086664398),
    @("FAPs","Inhbb","Acvr2a","FAPs",2,0.6666666666666666,0.3739756666666667,1.121927,0.1092711394000877,0.1092711394000878,3,1,27.61090666666666,82.83271999999999,0.4769721651858779,0.4769721651858778,10.32580722793778,92.93226505144,0.05211929195198773,0.05211929195198774),
    @("FAPs","Inhbb","Acvr2a","M2",  2,0.6666666666666666,0.3739756666666667,1.121927,0.1092711394000877,0.1092711394000878,3,1,4.423514,13.270542,0.07641520344774541,0.0764152034477454,1.654286597159333,14.888579374434,0.008349976348224653,0.008349976348224654),
    @("FAPs","Inhbb","Acvr2a","sCs", 2,0.6666666666666666,0.3739756666666667,1.121927,0.1092711394000877,0.1092711394000878,3,1,8.723210333333332,26.169631,0.1506914847198724,0.1506914847198724,3.262268399881889,29.360415598937,0.01646623023323137,0.01646623023323137),

    @("sCs","Inhbb","Acvr2a","ECs", 3,1,0.5233266666666666,1.56998,0.1529096843514326,0.1529096843514326,3,1,17.13024733333333,51.390742,0.2959211466465044,0.2959211466465043,8.964715236128887,80.68243712515999,0.04524920912663098,0.04524920912663098),
    @("sCs","Inhbb","Acvr2a","FAPs",3,1,0.5233266666666666,1.56998,0.1529096843514326,0.1529096843514326,3,1,27.61090666666666,82.83271999999999,0.4769721651858779,0.4769721651858778,14.44952374951111,130.0457137456,0.07293366322299195,0.07293366322299195),
    @("sCs","Inhbb","Acvr2a","M2",  3,1,0.5233266666666666,1.56998,0.1529096843514326,0.1529096843514326,3,1,4.423514,13.270542,0.07641520344774541,0.0764152034477454,2.314942836573333,20.83448552916,0.01168462463884525,0.01168462463884525),
    @("sCs","Inhbb","Acvr2a","sCs", 3,1,0.5233266666666666,1.56998,0.1529096843514326,0.1529096843514326,3,1,8.723210333333332,26.169631,0.1506914847198724,0.1506914847198724,4.565088586375555,41.08579727738,0.02304218736296442,0.02304218736296442)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$r = 2
foreach ($row in $rows) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value2 = $row[$i]
    }
    $r++
}
